$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98 (pushes existing rows 98:119 down to 99:120,
# carrying their formatting along automatically).
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new weekly price observation.
# Columns A,B,C,E,F,G,I,N,Q,R share the same market/category metadata as the
# surrounding rows (copy them straight from row 99, which now holds the data
# that used to live in row 98).
$ws.Cells.Item(98, 1).Value  = $ws.Cells.Item(99, 1).Value2   # A - Mercado ID
$ws.Cells.Item(98, 2).Value  = $ws.Cells.Item(99, 2).Value2   # B - Mercado
$ws.Cells.Item(98, 3).Value  = $ws.Cells.Item(99, 3).Value2   # C - Región
$ws.Cells.Item(98, 4).Value  = 45244                          # D - Fecha
$ws.Cells.Item(98, 5).Value  = $ws.Cells.Item(99, 5).Value2   # E - Codreg
$ws.Cells.Item(98, 6).Value  = $ws.Cells.Item(99, 6).Value2   # F - Categoría ID
$ws.Cells.Item(98, 7).Value  = $ws.Cells.Item(99, 7).Value2   # G - Categoría
$ws.Cells.Item(98, 8).Value  = "Sin especificar"              # H - Variedad
$ws.Cells.Item(98, 9).Value  = $ws.Cells.Item(99, 9).Value2   # I - Calidad
$ws.Cells.Item(98, 10).Value = 80                             # J - Volumen
$ws.Cells.Item(98, 11).Value = 25000                          # K - Precio mínimo
$ws.Cells.Item(98, 12).Value = 25000                          # L - Precio máximo
$ws.Cells.Item(98, 13).Value = 25000                          # M - Precio promedio ponderado
$ws.Cells.Item(98, 14).Value = $ws.Cells.Item(99, 14).Value2  # N - Unidad de comercialización
$ws.Cells.Item(98, 15).Value = "Región del Maule"              # O - Origen
$ws.Cells.Item(98, 16).Value = 1000                           # P - Precio $/Kg
$ws.Cells.Item(98, 17).Value = $ws.Cells.Item(99, 17).Value2  # Q - Kg o Unidades
$ws.Cells.Item(98, 18).Value = $ws.Cells.Item(99, 18).Value2  # R - Clasificación

# Keep the date cell's number format consistent with the rest of column D.
$ws.Cells.Item(98, 4).NumberFormat = $ws.Cells.Item(99, 4).NumberFormat
